$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (shifts old N/O/P -> O/P/Q)
$ws.Columns("N").Insert()

# The newly inserted column inherits its width from the column to its
# left (M) in real Excel's "insert column" behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet and select cell R6 on it,
# matching the new active-tab / selection state.
$ws.Select() | Out-Null
$ws.Range("R6").Select() | Out-Null
